$d = $word.ActiveDocument

# Update the "Updated:" date from 2023-07-29 to 2024-05-31 (the w:t run is
# wrapped by bookmarkStart/bookmarkEnd name="date", but Find/Replace on the
# document content operates on the visible text and leaves the bookmark
# anchors untouched).
$d.Content.Find.Execute("2023-07-29", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "2024-05-31", 2) | Out-Null

Write-Host "date updated"
